# Update the status ("E" column) of the action items for rows 52-56:
#  - Rows 52 & 53 (the two "Submit tickets ..." items) move from "Not Started" to "Complete"
#  - Rows 54, 55 & 56 (the "Upgrade ... tier." items) move from "Not Started" to "In Progress"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E52").Value = "Complete"
$ws.Range("E53").Value = "Complete"
$ws.Range("E54").Value = "In Progress"
$ws.Range("E55").Value = "In Progress"
$ws.Range("E56").Value = "In Progress"

# Reflect the final selection left on the sheet after the edits
$ws.Range("E55:E56").Select()
